# Apply updated crypto price/volume data (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.775.19'
$ws.Range('D3').Value = '2.577.82'
$ws.Range('E3').Value = '  -1.69%  '
$ws.Range('E4').Value = '  -0.19%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '585.22'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.67%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '168.96'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.24%  '
$ws.Range('E7').Value = '  -0.13%  '
$ws.Range('E8').Value = '  -0.79%  '
$ws.Range('D9').Value = '2.577.67'
$ws.Range('E9').Value = '  -1.76%  '
$ws.Range('E10').Value = '  +0.34%  '
$ws.Range('E11').Value = '  +0.38%  '
$ws.Range('E12').Value = '  -0.49%  '
$ws.Range('E13').Value = '  -0.86%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '26.83'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.76%  '
$ws.Range('D15').Value = '3.047.14'
$ws.Range('E15').Value = '  -2.01%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000179'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.29%  '
$ws.Range('D17').Value = '66.671.00'
$ws.Range('E17').Value = '  -0.36%  '
$ws.Range('D18').Value = '2.590.63'
$ws.Range('E18').Value = '  -1.31%  '
$ws.Range('E19').Value = '  -5.81%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.75'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.45%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '351.30'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.26%  '
$ws.Range('E22').Value = '  -1.64%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.61'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.88%  '
$ws.Range('E24').Value = '  -0.01%  '
$ws.Range('E25').Value = '  +0.77%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.95'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -8.60%  '
$ws.Range('D28').Value = '2.707.88'
$ws.Range('E28').Value = '  -1.95%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.13%  '
$ws.Range('D30').Value = '0.0₃0995'
$ws.Range('E30').Value = '  -0.83%  '
$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.31'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +5.12%  '
$ws.Range('B32').Value = 'Bittensor'
$ws.Range('C32').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '531.17'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.05%  '
$ws.Range('E33').Value = '  -1.48%  '
$ws.Range('E34').Value = '  -3.00%  '
$ws.Range('E35').Value = '  -2.86%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.999'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.09%  '
$ws.Range('E37').Value = '  -2.09%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '156.78'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.12%  '
$ws.Range('E39').Value = '  -1.12%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '18.33'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.19%  '
$ws.Range('E42').Value = '  -0.19%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.13'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.39%  '
$ws.Range('E44').Value = '  +0.03%  '
$ws.Range('E45').Value = '  +1.57%  '
$ws.Range('E46').Value = '  -2.91%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '149.24'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.31%  '
$ws.Range('E48').Value = '  -1.56%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.72'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.16%  '
$ws.Range('E50').Value = '  +1.29%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0763'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.96%  '
